$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cell A16 with the corrected timestamp value
$ws.Range("A16").Value = 45878.66688302084

# Append new row 17 with the latest reading
$ws.Range("A17").Value = 45878.70853278026
$ws.Range("B17").Value = 2025
$ws.Range("C17").Value = 37
$ws.Range("D17").Value = 17.78
$ws.Range("E17").Value = 80.23
$ws.Range("F17").Value = 138.87
$ws.Range("G17").Value = 12.74
$ws.Range("H17").Value = "ESE"
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = "17:00:17"

# Match the styling of the prior data row (column A carries the date/time style)
$ws.Range("A17").NumberFormat = $ws.Range("A16").NumberFormat
